# moovit_links.xlsx - "dan and Gur changed Dash2"
#
# Adds a second data row (route 11787, direction 2) whose link cell is a
# hyperlink to the Moovit route page, styled with the built-in "Hyperlink"
# cell style (blue/underlined).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of data under the existing header row.
$ws.Range("A2").Value = 11787
$ws.Range("B2").Value = 2

# Turn C2 into a hyperlink pointing at the Moovit line page. Excel.com
# interop auto-creates the shared string, the external relationship, the
# "Hyperlink" font/style (fontId 1 / xfId 1 / builtinId 8) and the
# <hyperlinks> entry on the sheet.
$null = $ws.Hyperlinks.Add($ws.Range("C2"), "https://moovitapp.com/israel-1/lines/787/204536/1088058/he?dayOffset=19536")

# Match the author's final cursor position / zoom level after the edit.
$null = $ws.Range("B3").Select()
$excel.ActiveWindow.Zoom = 251
